$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("R4").Value = 2021
$ws.Range("R5").Value = 47.8
$ws.Range("R6").Value = 20.7
$ws.Range("R7").Value = 9.8000000000000007
$ws.Range("R8").Value = 17.3

$ws.Range("Q4").Copy()
$ws.Range("R4").PasteSpecial(-4122)

$ws.Range("Q5").Copy()
$ws.Range("R5").PasteSpecial(-4122)

$ws.Range("Q6").Copy()
$ws.Range("R6").PasteSpecial(-4122)

$ws.Range("Q7").Copy()
$ws.Range("R7").PasteSpecial(-4122)

$ws.Range("Q8").Copy()
$ws.Range("R8").PasteSpecial(-4122)

$ws.Range("P10").Select()
